$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Runmode column values from "N" to "Y" for the Customer and Products rows
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "Y"

# Update the active selection to C4
$ws.Range("C4").Select()
